# #5: property boat&car done
# Turn the "汽車" (car) sheet's row 1 from a duplicate data row into a real
# header row, and add the common trailing metadata columns (H:N) -
# property_category / category / date / legislator_name / legislator_id /
# source_file / index - to the existing car/boat rows, matching the other
# property sheets (土地, 建物, 股票, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- Row 1: proper header labels (B1:N1) ----
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Match the bold/centered/bordered header style (style "1") used by the
# other header rows in this workbook.
$hdr = $ws.Range("B1:N1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# ---- Rows 2 & 3: append the trailing metadata columns (H:N) ----
$ws.Cells.Item(2, 8).Value = "land"
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(2, 10).Value = "2012-03-30"
$ws.Cells.Item(2, 11).Value = "黃志雄"
$ws.Cells.Item(2, 12).Value = 1366
$ws.Cells.Item(2, 13).Value = "tmpb8fa1"
$ws.Cells.Item(2, 14).Value = 31

$ws.Cells.Item(3, 8).Value = "land"
$ws.Cells.Item(3, 9).Value = "normal"
$ws.Cells.Item(3, 10).Value = "2012-03-30"
$ws.Cells.Item(3, 11).Value = "黃志雄"
$ws.Cells.Item(3, 12).Value = 1366
$ws.Cells.Item(3, 13).Value = "tmpb8fa1"
$ws.Cells.Item(3, 14).Value = 32
